# Timetabling GA result update: shift/re-assign a few "Kỹ năng mềm" and
# "Tiếng Anh chuyên ngành" sessions so the week no longer needs the
# T1 (17:30-19:30) CL10 / Ngô Văn I slot that previously lived in row 12,
# and instead folds that slot's class into the C2 (15:00-17:00) row that
# used to be merged with S2. Applies identically to every Tuan_N sheet.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    # Row 9: E9 picks up a new "Kỹ năng mềm" session (was blank); copy the
    # highlighted style from the cell that currently carries that same kind
    # of entry (G11) before we overwrite it.
    $ws.Range("G11").Copy()
    $ws.Range("E9").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("E9").Value2 = "Kỹ năng mềm`n(Lý thuyết)`nPhòng: R103`nGV: Hoàng Thị E"

    # I9: teacher reassigned.
    $ws.Range("I9").Value2 = "Tiếng Anh chuyên ngành`n(Lý thuyết)`nPhòng: R105`nGV: Hoàng Thị E"

    # A9:A10 was merged as a single "S2" slot label; split it and give A10
    # its own "C2 (15:00-17:00)" label, matching the A9/A11 slot-label style.
    $ws.Range("A9:A10").UnMerge()
    $ws.Range("A9").Copy()
    $ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("A10").Value2 = "C2`n(15:00-17:00)"

    # H10: room + teacher reassigned.
    $ws.Range("H10").Value2 = "Tiếng Anh chuyên ngành`n(Lý thuyết)`nPhòng: R103`nGV: Hoàng Thị E"

    # Row 11 becomes the old row-12 slot ("T1 17:30-19:30", class CL10).
    $ws.Range("A11").Value2 = "T1`n(17:30-19:30)"
    $ws.Range("B11").Value2 = "CL10"

    # The "Kỹ năng mềm" entry moves from G11 to H11 (new room R101).
    $ws.Range("D9").Copy()
    $ws.Range("G11").PasteSpecial(-4122)  # xlPasteFormats -> plain/blank style
    $ws.Range("G11").Value2 = ""

    $ws.Range("E9").Copy()
    $ws.Range("H11").PasteSpecial(-4122)  # xlPasteFormats -> highlighted style
    $ws.Range("H11").Value2 = "Kỹ năng mềm`n(Lý thuyết)`nPhòng: R101`nGV: Hoàng Thị E"

    # Old row 12 (T1 CL10 / Ngô Văn I / R104) is no longer needed as a
    # separate row now that its slot lives in row 11.
    $ws.Rows.Item(12).Delete()
}

Write-Host "Updated $($wb.Worksheets.Count) sheets"
